# Generate Report for Archive
# Two localization jobs (3f5e3040-...-9a8 and 96c2df2e-...-58a) moved from
# "Ready for handoff" to "In Translation" status. Update the Status column
# on the Overview roll-up sheet (both language columns) and on each
# per-language detail sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: rows 3 & 4 (zh-cn column B, de-de column C) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$overview.Range("B4").Value = $newStatus
$overview.Range("C4").Value = $newStatus

# --- zh-cn detail sheet: Status column (B) rows 3 & 4 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $newStatus
$zhcn.Range("B4").Value = $newStatus

# --- de-de detail sheet: Status column (B) rows 3 & 4 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $newStatus
$dede.Range("B4").Value = $newStatus
